# Atualizado por script em 04-01-2024 20:45
#
# This script:
#  1. Swaps the match-detail columns (F:V) between three pairs of rows
#     that were re-ordered by the scraper (12<->13, 96<->97, 101<->102).
#     Columns A:E (index/pais/torneio/temporada/data_partida) stay put.
#  2. Appends four newly scraped matches as rows 108-111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

# Scratch row far outside the used range, used as swap buffer.
$scratch = 500

function Swap-RowDetails($rowA, $rowB) {
    $rangeA = "F$($rowA):V$($rowA)"
    $rangeB = "F$($rowB):V$($rowB)"
    $rangeScratch = "F$($scratch):V$($scratch)"

    $ws.Range($rangeA).Copy()
    $ws.Range($rangeScratch).PasteSpecial($xlPasteValues)

    $ws.Range($rangeB).Copy()
    $ws.Range($rangeA).PasteSpecial($xlPasteValues)

    $ws.Range($rangeScratch).Copy()
    $ws.Range($rangeB).PasteSpecial($xlPasteValues)

    $ws.Range($rangeScratch).Clear()
}

Swap-RowDetails 12 13
Swap-RowDetails 96 97
Swap-RowDetails 101 102

$excel.CutCopyMode = $false

function Add-MatchRow(
    $row,
    $indice,
    $dataPartida,
    $home,
    $homeGols,
    $away,
    $awayGols,
    $homeOpenOdds,
    $homeOpenDh,
    $homeCloseOdds,
    $homeCloseDh,
    $drawOpenOdds,
    $drawOpenDh,
    $drawCloseOdds,
    $drawCloseDh,
    $awayOpenOdds,
    $awayOpenDh,
    $awayCloseOdds,
    $awayCloseDh,
    $url
) {
    # Clone the row-1-above formatting (bold/border index cell, datetime cell)
    # so the new rows carry the same cell styles as the rest of the sheet.
    $ws.Range("A$($row - 1)").Copy()
    $ws.Range("A$row").PasteSpecial($xlPasteFormats)
    $ws.Range("E$($row - 1)").Copy()
    $ws.Range("E$row").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $ws.Range("A$row").Value = $indice
    $ws.Range("B$row").Value = "morocco"
    $ws.Range("C$row").Value = "botola-pro"
    $ws.Range("D$row").Value = "2023-2024"
    $ws.Range("E$row").Value = $dataPartida
    $ws.Range("F$row").Value = $home
    $ws.Range("G$row").Value = $homeGols
    $ws.Range("H$row").Value = $away
    $ws.Range("I$row").Value = $awayGols
    $ws.Range("J$row").Value = $homeOpenOdds
    $ws.Range("K$row").Value = $homeOpenDh
    $ws.Range("L$row").Value = $homeCloseOdds
    $ws.Range("M$row").Value = $homeCloseDh
    $ws.Range("N$row").Value = $drawOpenOdds
    $ws.Range("O$row").Value = $drawOpenDh
    $ws.Range("P$row").Value = $drawCloseOdds
    $ws.Range("Q$row").Value = $drawCloseDh
    $ws.Range("R$row").Value = $awayOpenOdds
    $ws.Range("S$row").Value = $awayOpenDh
    $ws.Range("T$row").Value = $awayCloseOdds
    $ws.Range("U$row").Value = $awayCloseDh
    $ws.Range("V$row").Value = $url
}

Add-MatchRow 108 107 45295.66666666666 `
    "Union Touarga" 1 "Olympique de Safi" 0 `
    2.81 "02/01/2024 16:12" 2.75 "04/01/2024 15:56" `
    2.71 "02/01/2024 16:12" 2.96 "04/01/2024 15:54" `
    2.63 "02/01/2024 16:12" 2.7 "04/01/2024 15:56" `
    "https://www.betexplorer.com/football/morocco/botola-pro/union-touarga-olympique-de-safi/xrulz8ts/"

Add-MatchRow 109 108 45295.75 `
    "Berkane" 3 "Jeunesse Sportive Soualem" 1 `
    1.49 "02/01/2024 18:12" 1.38 "04/01/2024 17:58" `
    3.63 "02/01/2024 18:12" 4.28 "04/01/2024 17:58" `
    6.34 "02/01/2024 18:12" 9.12 "04/01/2024 17:58" `
    "https://www.betexplorer.com/football/morocco/botola-pro/berkane-jeunesse-sportive-soualem/AewdYoBg/"

Add-MatchRow 110 109 45295.75 `
    "Youssoufia Berrechid" 1 "Hassania Agadir" 2 `
    2.5 "02/01/2024 18:12" 2.54 "04/01/2024 17:59" `
    2.73 "02/01/2024 18:12" 2.46 "04/01/2024 17:58" `
    2.94 "02/01/2024 18:12" 3.69 "04/01/2024 17:59" `
    "https://www.betexplorer.com/football/morocco/botola-pro/youssoufia-berrechid-hassania-agadir/GYwhZRem/"

Add-MatchRow 111 110 45295.83333333334 `
    "FAR Rabat" 3 "Moghreb Tetouan" 0 `
    1.38 "02/01/2024 20:12" 1.28 "04/01/2024 19:53" `
    4.04 "02/01/2024 20:12" 5 "04/01/2024 19:59" `
    7.37 "02/01/2024 20:12" 11.3 "04/01/2024 19:59" `
    "https://www.betexplorer.com/football/morocco/botola-pro/far-rabat-moghreb-tetouan/MVt0X5Qa/"

$excel.CutCopyMode = $false
